# Updated cryptos list on Fri Oct 20 09:38:46 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures, and swaps the
# Kaspa/Aave rows (45/46) to reflect the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a literal text value (mirrors the original
    # inline-string cells) instead of letting Excel auto-convert numeric-
    # looking strings (e.g. "26.71") into real numbers.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.613.71"
$ws.Range("E2").Value = "  +4.34%  "
Set-TextValue $ws.Range("D3") "1.603.25"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("E6").Value = "  +7.69%  "
$ws.Range("E7").Value = "  -0.49%  "
Set-TextValue $ws.Range("D8") "26.71"
$ws.Range("E8").Value = "  +11.70%  "
Set-TextValue $ws.Range("D9") "0.250"
$ws.Range("E9").Value = "  +3.34%  "
Set-TextValue $ws.Range("D10") "0.0599"
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("E11").Value = "  +2.99%  "
Set-TextValue $ws.Range("D12") "1.833.48"
$ws.Range("E12").Value = "  +3.45%  "
Set-TextValue $ws.Range("D13") "1.595.15"
$ws.Range("E13").Value = "  +2.88%  "
Set-TextValue $ws.Range("D14") "29.646.59"
$ws.Range("E14").Value = "  +4.62%  "
Set-TextValue $ws.Range("D15") "3.76"
$ws.Range("E15").Value = "  +3.85%  "
$ws.Range("E16").Value = "  +3.60%  "
Set-TextValue $ws.Range("D17") "245.36"
$ws.Range("E17").Value = "  +7.50%  "
Set-TextValue $ws.Range("D18") "63.59"
$ws.Range("E18").Value = "  +4.50%  "
Set-TextValue $ws.Range("D19") "7.60"
$ws.Range("E19").Value = "  +3.50%  "
Set-TextValue $ws.Range("D20") "0.0₃0695"
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("E21").Value = "  -0.54%  "
Set-TextValue $ws.Range("D22") "4.05"
$ws.Range("E22").Value = "  +3.87%  "
Set-TextValue $ws.Range("D23") "9.27"
$ws.Range("E23").Value = "  +3.92%  "
Set-TextValue $ws.Range("D24") "2.11"
$ws.Range("E24").Value = "  +4.21%  "
Set-TextValue $ws.Range("D25") "155.94"
$ws.Range("E25").Value = "  +3.04%  "
Set-TextValue $ws.Range("D26") "15.36"
$ws.Range("E26").Value = "  +4.22%  "
$ws.Range("E27").Value = "  +5.85%  "
Set-TextValue $ws.Range("D28") "6.40"
$ws.Range("E28").Value = "  +2.58%  "
Set-TextValue $ws.Range("D29") "0.996"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E32").Value = "  +2.51%  "
Set-TextValue $ws.Range("D33") "1.438.95"
$ws.Range("E33").Value = "  +3.92%  "
Set-TextValue $ws.Range("D34") "3.11"
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("E36").Value = "  +10.85%  "
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("E39").Value = "  +2.53%  "
Set-TextValue $ws.Range("D40") "0.533"
$ws.Range("E40").Value = "  +4.46%  "
Set-TextValue $ws.Range("D41") "55.14"
$ws.Range("E41").Value = "  +27.54%  "
$ws.Range("E42").Value = "  +1.84%  "
Set-TextValue $ws.Range("D43") "0.800"
$ws.Range("E43").Value = "  +3.58%  "
$ws.Range("E44").Value = "  -0.56%  "

# Rows 45 and 46 swap places: Aave overtakes Kaspa in the ranking.
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "66.88"
$ws.Range("E45").Value = "  +7.91%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D46") "0.0467"
$ws.Range("E46").Value = "  +2.36%  "

Set-TextValue $ws.Range("D47") "5.32"
$ws.Range("E47").Value = "  -0.66%  "
Set-TextValue $ws.Range("D48") "1.743.98"
$ws.Range("E48").Value = "  +3.58%  "
Set-TextValue $ws.Range("D49") "86.28"
$ws.Range("E49").Value = "  +0.66%  "
Set-TextValue $ws.Range("D50") "0.836"
$ws.Range("E50").Value = "  -4.08%  "
Set-TextValue $ws.Range("D51") "0.0520"
$ws.Range("E51").Value = "  +1.80%  "
